# New crime data collected — update weekly CompStat figures (104th Precinct)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: report number and date range (shared, rich-text strings) ---
$ws.Range("A8").Text  = "Volume 29   Number  51"
$ws.Range("C9").Text  = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Row 15 (Rape) ---
$ws.Range("C15").Text  = "0"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 42.857142857142
$ws.Range("M15").Value = 5.263157894736
$ws.Range("N15").Value = -4.761904761904

# --- Row 16 ---
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 73.333333333333
$ws.Range("I16").Value = 217
$ws.Range("J16").Value = 141
$ws.Range("K16").Value = 53.900709219858
$ws.Range("L16").Value = 21.910112359550
$ws.Range("M16").Value = -17.175572519084
$ws.Range("N16").Value = -78.599605522682

# --- Row 17 ---
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 255
$ws.Range("K17").Value = 19.215686274509
$ws.Range("L17").Value = 6.293706293706
$ws.Range("M17").Value = 38.812785388127
$ws.Range("N17").Value = 0

# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 267
$ws.Range("J18").Value = 240
$ws.Range("K18").Value = 11.25
$ws.Range("L18").Value = -16.300940438871
$ws.Range("M18").Value = -40
$ws.Range("N18").Value = -86.187273667873

# --- Row 19 ---
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -47.368421052631
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 87
$ws.Range("H19").Value = -50.574712643678
$ws.Range("I19").Value = 629
$ws.Range("J19").Value = 577
$ws.Range("K19").Value = 9.012131715771
$ws.Range("L19").Value = 5.892255892255
$ws.Range("M19").Value = 48
$ws.Range("N19").Value = -5.838323353293

# --- Row 20 ---
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 43.75
$ws.Range("I20").Value = 309
$ws.Range("J20").Value = 210
$ws.Range("K20").Value = 47.142857142857
$ws.Range("L20").Value = 72.625698324022
$ws.Range("M20").Value = -13.687150837988
$ws.Range("N20").Value = -90.839015713015

# --- Row 21 (Total violent felony) ---
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 173
$ws.Range("H21").Value = -18.497109826589
$ws.Range("I21").Value = 1748
$ws.Range("J21").Value = 1441
$ws.Range("K21").Value = 21.304649548924
$ws.Range("L21").Value = 11.337579617834
$ws.Range("M21").Value = 0.923787528868
$ws.Range("N21").Value = -76.159301691216

# --- Row 22 (Robbery) ---
$ws.Range("C22").Value = 1
$ws.Range("D22").Text  = "0"
$ws.Range("E22").Text  = "***.*"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 88.888888888888
$ws.Range("L22").Value = -15
$ws.Range("M22").Value = -10.526315789473

# --- Row 24 ---
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -25.714285714285
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -4.225352112676
$ws.Range("I24").Value = 1479
$ws.Range("J24").Value = 1380
$ws.Range("K24").Value = 7.173913043478
$ws.Range("L24").Value = 1.440329218107
$ws.Range("M24").Value = 40.589353612167

# --- Row 25 ---
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -19.444444444444
$ws.Range("I25").Value = 523
$ws.Range("J25").Value = 489
$ws.Range("K25").Value = 6.952965235173
$ws.Range("L25").Value = 14.192139737991
$ws.Range("M25").Value = -23.088235294117

# --- Row 26 ---
$ws.Range("C26").Text  = "0"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 34
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 54.545454545454
$ws.Range("L26").Value = 25.925925925925

# --- Row 27 ---
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0

# --- Row 30 ---
$ws.Range("G30").Text = "0"
$ws.Range("H30").Text = "***.*"
